# Add files via upload
# - Rename the worksheet from "Sheet1" to "Conspiracy-Benchmark"
# - Move the active selection/active cell to G451 (the saved view's
#   topLeftCell for the frozen pane resets to A2 as part of this)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Conspiracy-Benchmark"

$ws.Range("G451").Select()
